$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2025-09-13 Saturday"

# Update table cell values
$t = $d.Tables.Item(1)
$values = @(
    "13+82=95", "36-18=18", "53+45=98", "90-18=72", "66-45=21",
    "64+28=92", "70-51=19", "95+0=95", "34-11=23", "13+77=90",
    "48+50=98", "94-53=41", "92-55=37", "52+17=69", "27+34=61",
    "31+23=54", "34+28=62", "50+0=50", "80-50=30", "97-12=85",
    "82-10=72", "53+11=64", "77+9=86", "20+57=77", "41-1=40",
    "39-28=11", "61-20=41", "67-56=11", "63-38=25", "27-2=25",
    "27-18=9", "21+36=57", "2+33=35", "64-44=20", "21+72=93",
    "16+22=38", "44-28=16", "59-34=25", "85+1=86", "89-86=3",
    "64+25=89", "53-5=48", "48+39=87", "16-15=1", "26+33=59",
    "22-4=18", "50+44=94", "65-48=17", "43+12=55", "86+13=99",
    "8+48=56", "50+19=69", "45-36=9", "74-35=39", "69-8=61",
    "23+55=78", "23+16=39", "17+57=74", "83-23=60", "21-17=4",
    "45+48=93", "51+47=98", "33+37=70", "53-24=29", "82+7=89",
    "87-64=23", "54-25=29", "88-18=70", "18-1=17", "68-37=31",
    "37+16=53", "26+34=60", "92-78=14", "74-40=34", "98-79=19",
    "31-19=12", "29+62=91", "62-49=13", "46+38=84", "88-18=70",
    "4+23=27", "10+37=47", "73-9=64", "4+67=71", "54-54=0",
    "14+79=93", "77-44=33", "10+41=51", "28+11=39", "82-37=45",
    "51-26=25", "45+11=56", "70+19=89", "10+18=28", "43-15=28",
    "24-0=24", "25+57=82", "94-63=31", "46-44=2", "5+85=90"
)

$numCols = 5
$numRows = 20
for ($r = 1; $r -le $numRows; $r++) {
    for ($c = 1; $c -le $numCols; $c++) {
        $idx = (($r - 1) * $numCols) + ($c - 1)
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
    }
}

Write-Output "Done updating document."